$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-13) holds a "Förändrad" (last changed) date that is
# stored as a date serial value. The automatic update bumps that date by
# one day (serial 46061 -> 46062, i.e. 2026-02-08 -> 2026-02-09) for every
# data row currently in the sheet.
for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value -ne $null) {
        $cell.Value = 46062
    }
}
